$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data values for rows 2-21 (columns A-T), derived from the target diff.
# Sending cluster / Ligand / Receptor / Target cluster / ... metrics
$data = @(
    @("ECs","Inhbb","Acvr2a","ECs",3,1,4.308396,12.925188,0.8865596696737007,0.8865596696737006,3,1,10.50827633333333,31.524829,0.1682660991018133,0.1682660991018134,45.273815721428,407.464341492852,0.1491779372369858,0.1491779372369858),
    @("ECs","Inhbb","Acvr2a","FAPs",3,1,4.308396,12.925188,0.8865596696737007,0.8865596696737006,3,1,30.94593233333333,92.837797,0.4955285863849104,0.4955285863849105,133.327331081204,1199.945979730836,0.4393156598592821,0.4393156598592821),
    @("ECs","Inhbb","Acvr2a","Inflammatory-Mac",3,1,4.308396,12.925188,0.8865596696737007,0.8865596696737006,3,1,6.495209666666667,19.485629,0.1040059814559238,0.1040059814559238,27.983935347028,251.855418123252,0.09220750856365283,0.09220750856365281),
    @("ECs","Inhbb","Acvr2a","MuSCs",3,1,4.308396,12.925188,0.8865596696737007,0.8865596696737006,3,1,9.909791666666667,29.729375,0.1586827309986352,0.1586827309986352,42.6953067775,384.2577609975,0.1406817095770707,0.1406817095770707),
    @("ECs","Inhbb","Acvr2a","Resolving-Mac",3,1,4.308396,12.925188,0.8865596696737007,0.8865596696737006,3,1,4.591137333333333,13.773412,0.07351660205871713,0.07351660205871713,19.780437722384,178.023939501456,0.06517685443670916,0.06517685443670916),
    @("FAPs","Inhbb","Acvr2a","ECs",1,0.3333333333333333,0.062595,0.187785,0.01288047861042144,0.01288047861042144,3,1,10.50827633333333,31.524829,0.1682660991018133,0.1682660991018134,0.657765557085,5.919890013765,0.002167347890339961,0.002167347890339962),
    @("FAPs","Inhbb","Acvr2a","FAPs",1,0.3333333333333333,0.062595,0.187785,0.01288047861042144,0.01288047861042144,3,1,30.94593233333333,92.837797,0.4955285863849104,0.4955285863849105,1.937060634405,17.433545709645,0.006382645357783212,0.006382645357783212),
    @("FAPs","Inhbb","Acvr2a","Inflammatory-Mac",1,0.3333333333333333,0.062595,0.187785,0.01288047861042144,0.01288047861042144,3,1,6.495209666666667,19.485629,0.1040059814559238,0.1040059814559238,0.406567649085,3.659108841765,0.001339646819498915,0.001339646819498915),
    @("FAPs","Inhbb","Acvr2a","MuSCs",1,0.3333333333333333,0.062595,0.187785,0.01288047861042144,0.01288047861042144,3,1,9.909791666666667,29.729375,0.1586827309986352,0.1586827309986352,0.620303409375,5.582730684375,0.00204390952247118,0.00204390952247118),
    @("FAPs","Inhbb","Acvr2a","Resolving-Mac",1,0.3333333333333333,0.062595,0.187785,0.01288047861042144,0.01288047861042144,3,1,4.591137333333333,13.773412,0.07351660205871713,0.07351660205871713,0.28738224138,2.58644017242,0.0009469290203281709,0.0009469290203281708),
    @("Inflammatory-Mac","Inhbb","Acvr2a","ECs",1,0.3333333333333333,0.05725033333333333,0.171751,0.01178068046871951,0.01178068046871951,3,1,10.50827633333333,31.524829,0.1682660991018133,0.1682660991018134,0.601602322842111,5.414420905579,0.001982289147236354,0.001982289147236354),
    @("Inflammatory-Mac","Inhbb","Acvr2a","FAPs",1,0.3333333333333333,0.05725033333333333,0.171751,0.01178068046871951,0.01178068046871951,3,1,30.94593233333333,92.837797,0.4955285863849104,0.4955285863849105,1.771664941394111,15.944984472547,0.005837663939316902,0.005837663939316901),
    @("Inflammatory-Mac","Inhbb","Acvr2a","Inflammatory-Mac",1,0.3333333333333333,0.05725033333333333,0.171751,0.01178068046871951,0.01178068046871951,3,1,6.495209666666667,19.485629,0.1040059814559238,0.1040059814559238,0.3718529184865555,3.346676266379,0.001225261234367805,0.001225261234367805),
    @("Inflammatory-Mac","Inhbb","Acvr2a","MuSCs",1,0.3333333333333333,0.05725033333333333,0.171751,0.01178068046871951,0.01178068046871951,3,1,9.909791666666667,29.729375,0.1586827309986352,0.1586827309986352,0.5673388761805556,5.106049885625,0.001869390549798693,0.001869390549798693),
    @("Inflammatory-Mac","Inhbb","Acvr2a","Resolving-Mac",1,0.3333333333333333,0.05725033333333333,0.171751,0.01178068046871951,0.01178068046871951,3,1,4.591137333333333,13.773412,0.07351660205871713,0.07351660205871713,0.2628441427124444,2.365597284412,0.0008660755979997532,0.0008660755979997531),
    @("MuSCs","Inhbb","Acvr2a","ECs",3,1,0.4314383333333334,1.294315,0.08877917124715834,0.08877917124715833,3,1,10.50827633333333,31.524829,0.1682660991018133,0.1682660991018134,4.533673227459444,40.803059047135,0.0149385248272512,0.0149385248272512),
    @("MuSCs","Inhbb","Acvr2a","FAPs",3,1,0.4314383333333334,1.294315,0.08877917124715834,0.08877917124715833,3,1,30.94593233333333,92.837797,0.4955285863849104,0.4955285863849105,13.35126146933944,120.161353224055,0.04399261722852826,0.04399261722852826),
    @("MuSCs","Inhbb","Acvr2a","Inflammatory-Mac",3,1,0.4314383333333334,1.294315,0.08877917124715834,0.08877917124715833,3,1,6.495209666666667,19.485629,0.1040059814559238,0.1040059814559238,2.802282433237222,25.220541899135,0.009233564838404233,0.009233564838404231),
    @("MuSCs","Inhbb","Acvr2a","MuSCs",3,1,0.4314383333333334,1.294315,0.08877917124715834,0.08877917124715833,3,1,9.909791666666667,29.729375,0.1586827309986352,0.1586827309986352,4.275464000347223,38.479176003125,0.01408772134929459,0.01408772134929459),
    @("MuSCs","Inhbb","Acvr2a","Resolving-Mac",3,1,0.4314383333333334,1.294315,0.08877917124715834,0.08877917124715833,3,1,4.591137333333333,13.773412,0.07351660205871713,0.07351660205871713,1.980792639197778,17.82713375278,0.006526743003680041,0.00652674300368004)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $row[$c - 1]
    }
}
